# Add a new row of test-case data to Table1 on the "Test Cases" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$table = $ws.ListObjects.Item("Table1")
$newRow = $table.ListRows.Add()

$ws.Range("A7").Value = "Pravash Test"
$ws.Range("B7").Value = "Add Element"
$ws.Range("C7").Value = "Add Element"
$ws.Range("D7").Value = "Next Step"

# "How to simulate?" column wraps text, matching the rest of that column.
$ws.Range("C7").WrapText = $true

[void]$ws.Range("A8").Select()
